$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9   = "c76eb7bc1a2e6e67034fc7d750762de9"
    17  = "6d8ffd1d66c53c710be8772851e1d28a"
    126 = "885125c12218fb55a9c17a473ea811f5"
    136 = "b7039b57dbda92005e340256ad999b90"
    159 = "e180276079263c04640119ac2f9a2356"
    169 = "57c8ebb0b1bfe05484cfbeee6e543676"
    183 = "1566ad624c9b683444f8640e7090cefd"
    200 = "d5ef55e36803ff9c65c83cdd13fffe52"
    228 = "fe38701a3da4b84079059572acfcc9b3"
    281 = "5303e7c7c414586e96e97fca9adc5a1a"
    302 = "c72ab92478c61d71a94c691b800f69f1"
    339 = "099ff95134ac2a6dda1c6112387b1c53"
    420 = "bf3569543f5afe0bd329968445d710df"
    464 = "88ca15026fa327f90edcf2607339c165"
    506 = "d11e0cd41977733cd7b40226af342944"
    507 = "28b7f4082aa807fa960d3091d6953006"
    524 = "b442e64966200cb4be835787721f9bae"
    558 = "cb211322d39ea5dcae043e1ec1002c9b"
    580 = "521ce29e8304ca26acab34907e3d08da"
    624 = "8eed330081db7ea415c2ac50c2458014"
    635 = "d450c3da6f90944d2dbd85eeeee6c17e"
    688 = "87f7d8c8d5f14748512c9245c79f6ea6"
    693 = "e992428de39ad6cc52cb72f089587295"
    708 = "c73244e4d02da93b2f5418460dd36c9d"
    711 = "d174fa8fbca0c777f41402c2571309ad"
    723 = "3d55dde6eea0e77c61e852a4347905de"
    827 = "6f14a86add7ba4c658e6672d743c2b75"
    838 = "10e0d3fcba82c94ccc94802d6c5c9179"
    843 = "08ec81e9257330f99b6ec686fc7b6d56"
    870 = "f5701873862730c7a2f060b7921941ab"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
